$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new values would otherwise be
# auto-parsed as numbers by Excel (column D is inline text in the source).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated values exactly as in the diff.
$ws.Range("D2").Value = '26.701.77'
$ws.Range("E2").Value = '  -6.87%  '
$ws.Range("D3").Value = '1.692.21'
$ws.Range("E3").Value = '  -6.13%  '
$ws.Range("D4").Value = '1.007'
$ws.Range("E4").Value = '  +0.36%  '
$ws.Range("D5").Value = '217.09'
$ws.Range("E5").Value = '  -6.08%  '
$ws.Range("B6").Value = 'USDC'
$ws.Range("C6").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D6").Value = '1.007'
$ws.Range("E6").Value = '  +0.29%  '
$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").Value = '0.4952'
$ws.Range("E7").Value = '  -16.65%  '
$ws.Range("D8").Value = '0.2587'
$ws.Range("E8").Value = '  -6.89%  '
$ws.Range("D9").Value = '21.72'
$ws.Range("E9").Value = '  -7.02%  '
$ws.Range("D10").Value = '0.06087'
$ws.Range("E10").Value = '  -11.06%  '
$ws.Range("D11").Value = '0.07297'
$ws.Range("E11").Value = '  -3.11%  '
$ws.Range("D12").Value = '1.685.17'
$ws.Range("E12").Value = '  -6.65%  '
$ws.Range("D13").Value = '4.428'
$ws.Range("D14").Value = '1.922.81'
$ws.Range("E14").Value = '  -6.10%  '
$ws.Range("D15").Value = '0.5703'
$ws.Range("E15").Value = '  -8.90%  '
$ws.Range("D16").Value = '0.000008173'
$ws.Range("E16").Value = '  -11.53%  '
$ws.Range("D17").Value = '64.75'
$ws.Range("E17").Value = '  -14.06%  '
$ws.Range("D18").Value = '26.706.30'
$ws.Range("E18").Value = '  -6.77%  '
$ws.Range("D19").Value = '5.001'
$ws.Range("E19").Value = '  -8.49%  '
$ws.Range("D20").Value = '1.007'
$ws.Range("E20").Value = '  +0.31%  '
$ws.Range("E21").Value = '  -6.20%  '
$ws.Range("D22").Value = '183.28'
$ws.Range("E22").Value = '  -13.07%  '
$ws.Range("D23").Value = '6.180'
$ws.Range("E23").Value = '  -9.73%  '
$ws.Range("D24").Value = '1.008'
$ws.Range("E24").Value = '  +0.36%  '
$ws.Range("D25").Value = '145.20'
$ws.Range("E25").Value = '  -5.98%  '
$ws.Range("D26").Value = '7.551'
$ws.Range("E26").Value = '  -3.62%  '
$ws.Range("D27").Value = '0.1130'
$ws.Range("E27").Value = '  -11.35%  '
$ws.Range("D28").Value = '15.21'
$ws.Range("E28").Value = '  -7.23%  '
$ws.Range("D29").Value = '1.316'
$ws.Range("E29").Value = '  -9.10%  '
$ws.Range("D30").Value = '0.05545'
$ws.Range("E30").Value = '  -11.55%  '
$ws.Range("D31").Value = '1.323'
$ws.Range("E31").Value = '  -6.82%  '
$ws.Range("D32").Value = '3.458'
$ws.Range("E32").Value = '  -8.23%  '
$ws.Range("D33").Value = '3.450'
$ws.Range("E33").Value = '  -7.61%  '
$ws.Range("D34").Value = '1.652'
$ws.Range("E34").Value = '  -3.74%  '
$ws.Range("D35").Value = '1.005'
$ws.Range("E35").Value = '  -4.73%  '
$ws.Range("D36").Value = '2.408'
$ws.Range("E36").Value = '  -3.67%  '
$ws.Range("D37").Value = '0.5855'
$ws.Range("E37").Value = '  -8.51%  '
$ws.Range("D38").Value = '2.629'
$ws.Range("E38").Value = '  -3.29%  '
$ws.Range("D39").Value = '0.01577'
$ws.Range("E39").Value = '  -7.73%  '
$ws.Range("D40").Value = '1.066.93'
$ws.Range("E40").Value = '  -6.18%  '
$ws.Range("D41").Value = '5.865'
$ws.Range("E41").Value = '  -8.25%  '
$ws.Range("D42").Value = '0.8515'
$ws.Range("E42").Value = '  -1.59%  '
$ws.Range("D43").Value = '1.005'
$ws.Range("E43").Value = '  +0.07%  '
$ws.Range("D44").Value = '98.19'
$ws.Range("E44").Value = '  -2.39%  '
$ws.Range("D45").Value = '1.854.17'
$ws.Range("E45").Value = '  -5.45%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = '56.27'
$ws.Range("E46").Value = '  -6.93%  '
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").Value = '0.00000000106'
$ws.Range("E47").Value = '  -5.16%  '
$ws.Range("D48").Value = '1.003'
$ws.Range("E48").Value = '  -0.26%  '
$ws.Range("D49").Value = '8.049'
$ws.Range("E49").Value = '  -3.31%  '
$ws.Range("D50").Value = '0.4337'
$ws.Range("E50").Value = '  -3.67%  '
$ws.Range("D51").Value = '0.05201'
$ws.Range("E51").Value = '  -4.90%  '
